$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the date-column formatting (numFmt/font/border/alignment) from the
# last existing row onto the new row's date cell before writing values.
$ws.Range("A38").Copy()
$ws.Range("A39").PasteSpecial(-4122)  # xlPasteFormats

# Append the new forecast vector row (row 39)
$ws.Range("A39").Value = 45986
$ws.Range("B39").Value = 2025
$ws.Range("C39").Value = 0.3398512689293476
$ws.Range("D39").Value = 2026
$ws.Range("E39").Value = 0.8571438361188566
